$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 344, shifting existing rows 344-403 down to 345-404
$ws.Rows("344:344").Insert()

# Populate the newly inserted row 344 with the new weekly data record
$ws.Cells.Item(344, 1).Value2 = 4
$ws.Cells.Item(344, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(344, 3).Value2 = "Los Lagos"
$ws.Cells.Item(344, 4).Value2 = 45180
$ws.Cells.Item(344, 5).Value2 = 10
$ws.Cells.Item(344, 6).Value2 = 100112039
$ws.Cells.Item(344, 7).Value2 = "Ciboulette"
$ws.Cells.Item(344, 8).Value2 = "Sin especificar"
$ws.Cells.Item(344, 9).Value2 = "Primera"
$ws.Cells.Item(344, 10).Value2 = 80
$ws.Cells.Item(344, 11).Value2 = 3500
$ws.Cells.Item(344, 12).Value2 = 3500
$ws.Cells.Item(344, 13).Value2 = 3500
$ws.Cells.Item(344, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(344, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(344, 16).Value2 = 1167
$ws.Cells.Item(344, 17).Value2 = 3
$ws.Cells.Item(344, 18).Value2 = "Hortaliza"
